$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") '97.599.27'
$ws.Range("E2").Value = '  +4.45%  '

# Row 3
Set-TextValue $ws.Range("D3") '3.136.02'
$ws.Range("E3").Value = '  +0.61%  '

# Row 4
$ws.Range("E4").Value = '  -0.08%  '

# Row 5
Set-TextValue $ws.Range("D5") '242.44'
$ws.Range("E5").Value = '  +1.97%  '

# Row 6
Set-TextValue $ws.Range("D6") '609.56'
$ws.Range("E6").Value = '  -0.89%  '

# Row 7
$ws.Range("E7").Value = '  +0.64%  '

# Row 8
Set-TextValue $ws.Range("D8") '0.384'
$ws.Range("E8").Value = '  -2.19%  '

# Row 9
$ws.Range("E9").Value = '  +0.03%  '

# Row 10
Set-TextValue $ws.Range("D10") '3.134.10'
$ws.Range("E10").Value = '  +0.64%  '

# Row 11
Set-TextValue $ws.Range("D11") '0.782'
$ws.Range("E11").Value = '  -7.19%  '

# Row 12
$ws.Range("E12").Value = '  -0.26%  '

# Row 13
Set-TextValue $ws.Range("D13") '97.215.03'
$ws.Range("E13").Value = '  +4.37%  '

# Row 14
Set-TextValue $ws.Range("D14") '0.0000239'
$ws.Range("E14").Value = '  -2.61%  '

# Row 15
$ws.Range("E15").Value = '  +0.08%  '

# Row 16
Set-TextValue $ws.Range("D16") '33.75'
$ws.Range("E16").Value = '  -4.33%  '

# Row 17
Set-TextValue $ws.Range("D17") '3.717.03'
$ws.Range("E17").Value = '  +0.65%  '

# Row 18
Set-TextValue $ws.Range("D18") '3.129.63'
$ws.Range("E18").Value = '  +0.87%  '

# Row 19
Set-TextValue $ws.Range("D19") '518.50'
$ws.Range("E19").Value = '  +16.97%  '

# Row 20
Set-TextValue $ws.Range("D20") '3.40'
$ws.Range("E20").Value = '  -10.47%  '

# Row 21
Set-TextValue $ws.Range("D21") '14.38'
$ws.Range("E21").Value = '  -3.40%  '

# Row 22
Set-TextValue $ws.Range("D22") '5.62'
$ws.Range("E22").Value = '  -7.25%  '

# Row 23
$ws.Range("E23").Value = '  -5.15%  '

# Row 24
$ws.Range("E24").Value = '  -4.41%  '

# Row 25
Set-TextValue $ws.Range("D25") '88.27'
$ws.Range("E25").Value = '  +2.86%  '

# Row 26
$ws.Range("E26").Value = '  -6.43%  '

# Row 27
$ws.Range("E27").Value = '  -11.29%  '

# Row 28
Set-TextValue $ws.Range("D28") '3.290.85'
$ws.Range("E28").Value = '  +0.26%  '

# Row 29
$ws.Range("E29").Value = '  +0.15%  '

# Row 30
$ws.Range("E30").Value = '  -0.30%  '

# Row 31
$ws.Range("E31").Value = '  -4.02%  '

# Row 32
Set-TextValue $ws.Range("D32") '0.121'
$ws.Range("E32").Value = '  -3.03%  '

# Row 33
$ws.Range("E33").Value = '  -0.79%  '

# Row 34
$ws.Range("E34").Value = '  -4.30%  '

# Row 35
Set-TextValue $ws.Range("D35") '26.53'
$ws.Range("E35").Value = '  +2.10%  '

# Row 36
$ws.Range("E36").Value = '  -6.07%  '

# Row 37
Set-TextValue $ws.Range("D37") '7.23'
$ws.Range("E37").Value = '  -9.80%  '

# Row 38
$ws.Range("B38").Value = 'WhiteBITCoin'
$ws.Range("C38").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
Set-TextValue $ws.Range("D38") '24.32'
$ws.Range("E38").Value = '  +1.34%  '

# Row 39
$ws.Range("B39").Value = 'PancakeSwap'
$ws.Range("C39").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue $ws.Range("D39") '1.87'
$ws.Range("E39").Value = '  -2.29%  '

# Row 40
$ws.Range("B40").Value = 'Bittensor'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue $ws.Range("D40") '467.95'
$ws.Range("E40").Value = '  -2.01%  '

# Row 41
$ws.Range("B41").Value = 'PolygonEcosystemToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
Set-TextValue $ws.Range("D41") '0.432'
$ws.Range("E41").Value = '  -3.39%  '

# Row 42
$ws.Range("E42").Value = '  -6.85%  '

# Row 43
$ws.Range("E43").Value = '  -10.42%  '

# Row 44
$ws.Range("E44").Value = '  +0.00%  '

# Row 45
$ws.Range("E45").Value = '  -6.97%  '

# Row 46
Set-TextValue $ws.Range("D46") '162.42'
$ws.Range("E46").Value = '  +2.12%  '

# Row 47
$ws.Range("E47").Value = '  -1.55%  '

# Row 48
$ws.Range("E48").Value = '  +0.99%  '

# Row 49
Set-TextValue $ws.Range("D49") '4.46'
$ws.Range("E49").Value = '  +0.81%  '

# Row 50
$ws.Range("E50").Value = '  +0.22%  '

# Row 51
$ws.Range("E51").Value = '  +0.03%  '
